# Updated cryptos list - applies new Price/Volume(1h) figures plus the
# row-25/26 Monero<->Stellar swap, matching the scraper run for this commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.371.10"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").Value = "1.870.89"
$ws.Range("E3").Value = "  +0.48%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9991"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7156"
$ws.Range("E5").Value = "  +2.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "238.44"
$ws.Range("E6").Value = "  +0.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9998"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07868"
$ws.Range("E8").Value = "  -4.33%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3073"
$ws.Range("E9").Value = "  +0.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.28"
$ws.Range("E10").Value = "  +8.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08183"
$ws.Range("E11").Value = "  +0.14%  "
$ws.Range("D12").Value = "1.865.00"
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.255"
$ws.Range("E13").Value = "  +1.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7210"
$ws.Range("E14").Value = "  +0.47%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.45"
$ws.Range("E15").Value = "  +0.07%  "
$ws.Range("D16").Value = "29.388.80"
$ws.Range("E16").Value = "  +0.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.831"
$ws.Range("E17").Value = "  +0.77%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "242.49"
$ws.Range("E18").Value = "  +2.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007820"
$ws.Range("E19").Value = "  -0.55%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.28"
$ws.Range("E20").Value = "  -0.91%  "
$ws.Range("D21").Value = "2.117.08"
$ws.Range("E21").Value = "  +0.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9997"
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9996"
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.630"
$ws.Range("E24").Value = "  +2.11%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "162.02"
$ws.Range("E25").Value = "  -0.23%  "
$ws.Range("B26").Value = "Stellar"
$ws.Range("C26").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1463"
$ws.Range("E26").Value = "  +1.67%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.954"
$ws.Range("E27").Value = "  -0.44%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.21"
$ws.Range("E28").Value = "  +0.44%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.926"
$ws.Range("E29").Value = "  -2.82%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.366"
$ws.Range("E30").Value = "  -4.47%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.334"
$ws.Range("E32").Value = "  -2.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.060"
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05234"
$ws.Range("E34").Value = "  +0.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.190"
$ws.Range("E35").Value = "  +1.52%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7206"
$ws.Range("E36").Value = "  +2.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.003"
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.673"
$ws.Range("E38").Value = "  +0.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01853"
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.706"
$ws.Range("E40").Value = "  -0.97%  "
$ws.Range("D41").Value = "1.179.77"
$ws.Range("E41").Value = "  +3.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9174"
$ws.Range("E42").Value = "  -0.37%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.017"
$ws.Range("E43").Value = "  +0.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4299"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "71.35"
$ws.Range("E45").Value = "  +0.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9999"
$ws.Range("E46").Value = "  -0.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.24"
$ws.Range("E47").Value = "  -0.57%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5359"
$ws.Range("E48").Value = "  -0.98%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.765"
$ws.Range("E49").Value = "  -0.56%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.228"
$ws.Range("E50").Value = "  +0.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.052"
$ws.Range("E51").Value = "  +1.09%  "
